$d = $word.ActiveDocument

# 1. Update the "5. Princípios SOLID" contributors line: reorder / add
#    "Gabriel Bezerra" to the list of names.
$find = $d.Content.Find
$replaced = $find.Execute("Millena Costa, Gabriel Ferreira e Danilo Murilo.", $true, $false, $false, $false, $false, $true, 1, $false, "Danilo Murilo, Millena Costa, Gabriel Bezerra e Gabriel Ferreira.", 2)
Write-Output ("Contributors line replaced: " + $replaced)

# 2. Append a new "7. Diagrama de Classes:" section at the very end of the
#    document (after the "Back-End: Gabriel Bezerra" line), made up of:
#      - a blank spacer paragraph
#      - a bold heading paragraph
#      - a blank spacer paragraph
#      - a paragraph naming the contributor, "Gabriel Bezerra,"
$endRange = $d.Content
$endRange.Collapse(0)

$p1 = '<w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
$p2 = '<w:p><w:pPr><w:rPr><w:b w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">7. Diagrama de Classes:</w:t></w:r></w:p>'
$p3 = '<w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
$p4 = '<w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Gabriel Bezerra,</w:t></w:r><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'

$newParasBody = $p1 + $p2 + $p3 + $p4

$xmlHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$xmlFragment = $xmlHeader + $newParasBody + $xmlFooter

$endRange.InsertXML($xmlFragment)

Write-Output "Edit complete"
